$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Update timestamp in title cell (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 27 de Mayo de 2020 a las 11:35"

# --- Swap the order of "Albania" and "Guinea Ecuatorial" entries ---
# Row 112 used to be "Guinea Ecuatorial" and row 113 used to be "Albania".
# After the update, row 112 becomes "Albania" (with refreshed stats) and
# row 113 becomes "Guinea Ecuatorial" (keeping its previous stats).
$ws.Range("A112").Value = "Albania"
$ws.Range("A113").Value = "Guinea Ecuatorial"

# --- Update numeric statistics for the affected country rows ---

# Row 61: Malasia
$ws.Range("B61").Value = 7619
$ws.Range("C61").Value = 15
$ws.Range("D61").Value = 6083
$ws.Range("E61").Value = 1421

# Row 100: Eslovenia
$ws.Range("B100").Value = 1471
$ws.Range("C100").Value = 2
$ws.Range("D100").Value = 1354
$ws.Range("E100").Value = 9

# Row 109: Hong Kong
$ws.Range("B109").Value = 1067
$ws.Range("C109").Value = 1
$ws.Range("D109").Value = 1034

# Row 112: Albania (new data)
$ws.Range("B112").Value = 1050
$ws.Range("C112").Value = 21
$ws.Range("D112").Value = 812
$ws.Range("E112").Value = 205
$ws.Range("H112").Value = 33

# Row 113: Guinea Ecuatorial (old data, unchanged values moved here)
$ws.Range("B113").Value = 1043
$ws.Range("C113").Value = 0
$ws.Range("D113").Value = 165
$ws.Range("E113").Value = 866
$ws.Range("H113").Value = 12

# Row 119: Burkina Faso
$ws.Range("B119").Value = 845
$ws.Range("C119").Value = 13
$ws.Range("E119").Value = 120
$ws.Range("G119").Value = 1
$ws.Range("H119").Value = 53

# Row 163: Brunei
$ws.Range("E163").Value = 2
$ws.Range("G163").Value = 1
$ws.Range("H163").Value = 2
